$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 2.5
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 126334034
$ws.Range("I2").Value = "REKR"

$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 2.5
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 126334034
$ws.Range("I3").Value = "REKR"

$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 2.5
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 126334034
$ws.Range("I4").Value = "REKR"

$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2.5
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 126334034
$ws.Range("I5").Value = "REKR"

$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 2.5
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 126334034
$ws.Range("I6").Value = "REKR"

$ws.Range("D7").Value = 2.019000053405762
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 5.5
$ws.Range("G7").Value = 1.509999990463257
$ws.Range("H7").Value = 126334034
$ws.Range("I7").Value = "REKR"

$ws.Range("D8").Value = 4.699999809265137
$ws.Range("E8").Value = 3.75
$ws.Range("F8").Value = 5.5
$ws.Range("G8").Value = 3.230000019073486
$ws.Range("H8").Value = 126334034
$ws.Range("I8").Value = "REKR"

$ws.Range("D9").Value = 3.400000095367432
$ws.Range("E9").Value = 2.926000118255615
$ws.Range("F9").Value = 3.901000022888184
$ws.Range("G9").Value = 2.299999952316284
$ws.Range("H9").Value = 126334034
$ws.Range("I9").Value = "REKR"

$ws.Range("D10").Value = 1.651000022888184
$ws.Range("E10").Value = 1.302999973297119
$ws.Range("F10").Value = 2.062999963760376
$ws.Range("G10").Value = 1.179999947547913
$ws.Range("H10").Value = 126334034
$ws.Range("I10").Value = "REKR"

$ws.Range("D11").Value = 1.269999980926514
$ws.Range("E11").Value = 0.7900000214576721
$ws.Range("F11").Value = 3.5
$ws.Range("G11").Value = 0.6809999942779541
$ws.Range("H11").Value = 126334034
$ws.Range("I11").Value = "REKR"

$ws.Range("D12").Value = 0.718999981880188
$ws.Range("E12").Value = 0.6669999957084656
$ws.Range("F12").Value = 0.7950000166893005
$ws.Range("G12").Value = 0.5099999904632568
$ws.Range("H12").Value = 126334034
$ws.Range("I12").Value = "REKR"

$ws.Range("D13").Value = 0.7039999961853027
$ws.Range("E13").Value = 0.7200000286102295
$ws.Range("F13").Value = 0.7200000286102295
$ws.Range("G13").Value = 0.5899999737739563
$ws.Range("H13").Value = 126334034
$ws.Range("I13").Value = "REKR"

$ws.Range("D14").Value = 1.889999985694885
$ws.Range("E14").Value = 4.320000171661377
$ws.Range("F14").Value = 5.440000057220459
$ws.Range("G14").Value = 1.409999966621399
$ws.Range("H14").Value = 126334034
$ws.Range("I14").Value = "REKR"

$ws.Range("D15").Value = 2.150000095367432
$ws.Range("E15").Value = 2.099999904632568
$ws.Range("F15").Value = 2.529999971389771
$ws.Range("G15").Value = 1.850000023841858
$ws.Range("H15").Value = 126334034
$ws.Range("I15").Value = "REKR"

$ws.Range("D16").Value = 3.75
$ws.Range("E16").Value = 4.619999885559082
$ws.Range("F16").Value = 4.829999923706055
$ws.Range("G16").Value = 3.230000019073486
$ws.Range("H16").Value = 126334034
$ws.Range("I16").Value = "REKR"

$ws.Range("D17").Value = 3.410000085830688
$ws.Range("E17").Value = 3.890000104904175
$ws.Range("F17").Value = 4.119999885559082
$ws.Range("G17").Value = 3.150000095367432
$ws.Range("H17").Value = 126334034
$ws.Range("I17").Value = "REKR"

$ws.Range("D18").Value = 4.019999980926514
$ws.Range("E18").Value = 3.970000028610229
$ws.Range("F18").Value = 4.570000171661377
$ws.Range("G18").Value = 3.420000076293945
$ws.Range("H18").Value = 126334034
$ws.Range("I18").Value = "REKR"

$ws.Range("D19").Value = 5.75
$ws.Range("E19").Value = 4.159999847412109
$ws.Range("F19").Value = 6.849999904632568
$ws.Range("G19").Value = 3.089999914169312
$ws.Range("H19").Value = 126334034
$ws.Range("I19").Value = "REKR"

$ws.Range("D20").Value = 8.079999923706055
$ws.Range("E20").Value = 11.89999961853027
$ws.Range("F20").Value = 13.69999980926514
$ws.Range("G20").Value = 7.389999866485596
$ws.Range("H20").Value = 126334034
$ws.Range("I20").Value = "REKR"

$ws.Range("D21").Value = 20.10000038146973
$ws.Range("E21").Value = 23.45000076293945
$ws.Range("F21").Value = 25.3799991607666
$ws.Range("G21").Value = 18.47999954223633
$ws.Range("H21").Value = 126334034
$ws.Range("I21").Value = "REKR"

$ws.Range("D22").Value = 10.9350004196167
$ws.Range("E22").Value = 7.880000114440918
$ws.Range("F22").Value = 10.9350004196167
$ws.Range("G22").Value = 7.199999809265137
$ws.Range("H22").Value = 126334034
$ws.Range("I22").Value = "REKR"

$ws.Range("D23").Value = 11.59000015258789
$ws.Range("E23").Value = 13.60000038146973
$ws.Range("F23").Value = 13.89999961853027
$ws.Range("G23").Value = 10.22999954223633
$ws.Range("H23").Value = 126334034
$ws.Range("I23").Value = "REKR"

$ws.Range("D24").Value = 6.599999904632568
$ws.Range("E24").Value = 4.650000095367432
$ws.Range("F24").Value = 7.329999923706055
$ws.Range("G24").Value = 3.75
$ws.Range("H24").Value = 126334034
$ws.Range("I24").Value = "REKR"

$ws.Range("D25").Value = 4.119999885559082
$ws.Range("E25").Value = 2.970000028610229
$ws.Range("F25").Value = 5.079999923706055
$ws.Range("G25").Value = 2.910000085830688
$ws.Range("H25").Value = 126334034
$ws.Range("I25").Value = "REKR"

$ws.Range("D26").Value = 1.759999990463257
$ws.Range("E26").Value = 1.809999942779541
$ws.Range("F26").Value = 2.25
$ws.Range("G26").Value = 1.539999961853027
$ws.Range("H26").Value = 126334034
$ws.Range("I26").Value = "REKR"

$ws.Range("D27").Value = 1.029999971389771
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 1.330000042915344
$ws.Range("G27").Value = 0.9110000133514404
$ws.Range("H27").Value = 126334034
$ws.Range("I27").Value = "REKR"

$ws.Range("D28").Value = 1.259999990463257
$ws.Range("E28").Value = 1.695000052452087
$ws.Range("F28").Value = 2.234999895095825
$ws.Range("G28").Value = 1.174999952316284
$ws.Range("H28").Value = 126334034
$ws.Range("I28").Value = "REKR"

$ws.Range("D29").Value = 1.25
$ws.Range("E29").Value = 1.210000038146973
$ws.Range("F29").Value = 1.320000052452087
$ws.Range("G29").Value = 0.9980000257492064
$ws.Range("H29").Value = 126334034
$ws.Range("I29").Value = "REKR"

$ws.Range("D30").Value = 1.809999942779541
$ws.Range("E30").Value = 3.039999961853028
$ws.Range("F30").Value = 3.5
$ws.Range("G30").Value = 1.809999942779541
$ws.Range("H30").Value = 126334034
$ws.Range("I30").Value = "REKR"

$ws.Range("D31").Value = 2.819999933242798
$ws.Range("E31").Value = 2.75
$ws.Range("F31").Value = 3.559999942779541
$ws.Range("G31").Value = 2.519999980926514
$ws.Range("H31").Value = 126334034
$ws.Range("I31").Value = "REKR"

$ws.Range("D32").Value = 3.289999961853028
$ws.Range("E32").Value = 3.150000095367432
$ws.Range("F32").Value = 3.920000076293945
$ws.Range("G32").Value = 2.960000038146973
$ws.Range("H32").Value = 126334034
$ws.Range("I32").Value = "REKR"

$ws.Range("D33").Value = 2.299999952316284
$ws.Range("E33").Value = 1.769999980926514
$ws.Range("F33").Value = 2.569999933242798
$ws.Range("G33").Value = 1.730000019073486
$ws.Range("H33").Value = 126334034
$ws.Range("I33").Value = "REKR"

$ws.Range("D34").Value = 1.549999952316284
$ws.Range("E34").Value = 1.820000052452088
$ws.Range("F34").Value = 2.069999933242798
$ws.Range("G34").Value = 1.490000009536743
$ws.Range("H34").Value = 126334034
$ws.Range("I34").Value = "REKR"

$ws.Range("D35").Value = 1.169999957084656
$ws.Range("E35").Value = 1.049999952316284
$ws.Range("F35").Value = 1.220000028610229
$ws.Range("G35").Value = 1.009999990463257
$ws.Range("H35").Value = 126334034
$ws.Range("I35").Value = "REKR"

$ws.Range("D36").Value = 1.690000057220459
$ws.Range("E36").Value = 2.029999971389771
$ws.Range("F36").Value = 2.585000038146973
$ws.Range("G36").Value = 1.399999976158142
$ws.Range("H36").Value = 126334034
$ws.Range("I36").Value = "REKR"

$ws.Range("D37").Value = 0.8809999823570251
$ws.Range("E37").Value = 1.049999952316284
$ws.Range("F37").Value = 1.080000042915344
$ws.Range("G37").Value = 0.6209999918937683
$ws.Range("H37").Value = 126334034
$ws.Range("I37").Value = "REKR"

$ws.Range("D38").Value = 1.159999966621399
$ws.Range("E38").Value = 1.129999995231628
$ws.Range("F38").Value = 1.419999957084656
$ws.Range("G38").Value = 1.080000042915344
$ws.Range("H38").Value = 126334034
$ws.Range("I38").Value = "REKR"

